# Updated cryptos list on Sun Nov 12 13:42:17 UTC 2023 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) figures for
# each coin row (rows 2-51) of the active worksheet to the latest snapshot.
#
# Price cells must stay plain TEXT (they already are, e.g. "37.127.59"),
# so writes briefly force a text NumberFormat to stop Excel re-interpreting
# numeric-looking strings (like "59.06") as real numbers, then restore the
# original General/Normal formatting so no extra style is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.102.94"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.30%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.046.64"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.75%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.63"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.61%  "

$ws.Range("E6").Value = "  -1.68%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.06"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.16%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("E9").Value = "  +1.42%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0786"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.98%  "

$ws.Range("E11").Value = "  +0.81%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.77"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.26%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.346.63"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.45%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.835"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.06%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.73"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.29%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.051.27"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "17.91"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +22.90%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.105.97"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.99"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0897"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.67%  "

$ws.Range("E21").Value = "  -0.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "237.06"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.91%  "

$ws.Range("E23").Value = "  -0.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.45"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.31%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "169.11"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.64%  "

$ws.Range("E26").Value = "  +7.91%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.37"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.89%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.04"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.90%  "

$ws.Range("E29").Value = "  -0.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.14"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.14%  "

$ws.Range("E31").Value = "  +2.88%  "

$ws.Range("E32").Value = "  -1.54%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.52"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.28%  "

$ws.Range("E34").Value = "  +1.58%  "

$ws.Range("E35").Value = "  -0.07%  "

$ws.Range("E36").Value = "  -2.99%  "

$ws.Range("E37").Value = "  -1.32%  "

$ws.Range("E38").Value = "  -3.85%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.33"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.65%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.17"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +13.35%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.14"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +14.85%  "

$ws.Range("E42").Value = "  -1.86%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.43"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.31%  "

$ws.Range("E44").Value = "  -1.64%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "96.07"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.44%  "

$ws.Range("E46").Value = "  -1.98%  "

$ws.Range("E47").Value = "  -0.33%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.282.83"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.67%  "

$ws.Range("E49").Value = "  -1.85%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.233.62"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.26%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.56"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -19.22%  "
